# [REF] Tools refactoring: z0lib
#
# Update the account_move sample data: fix the ref/fiscal year for the
# z0bug.move_2 / z0bug.move_3 entries and append new journal entries
# (z0bug.move_4 .. z0bug.move_10) for the 2020/2021 rent quarters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (z0bug.move_2): date/ref changed, fiscalyear_id changed ---
$ws.Range("C3").Value = "<1-01-15"
$ws.Range("D3").Value = "Pagata RA"
$ws.Range("F3").Value = "z0bug.fy_2020"

# --- Row 4 (z0bug.move_3): date/ref changed, fiscalyear_id changed ---
$ws.Range("C4").Value = "<1-01-31"
$ws.Range("D4").Value = "Affitto 1.mo Trimestre 2020"
$ws.Range("F4").Value = "z0bug.fy_2020"

# --- New rows 5..11 (z0bug.move_4 .. z0bug.move_10), same layout as the
#     existing ones: A id, B blank, C date, D ref, E journal_id,
#     F fiscalyear_id, G type ---
$rows = @(
    @{ Row = 5;  A = "z0bug.move_4";  C = "<1-04-05";   D = "Affitto 2.do Trimestre 2020"; F = "z0bug.fy_2020" },
    @{ Row = 6;  A = "z0bug.move_5";  C = "<1-07-05";   D = "Affitto 3.zo Trimestre 2020"; F = "z0bug.fy_2020" },
    @{ Row = 7;  A = "z0bug.move_6";  C = "<1-10-05";   D = "Affitto 4.to Trimestre 2020"; F = "z0bug.fy_2020" },
    @{ Row = 8;  A = "z0bug.move_7";  C = "####-01-31"; D = "Affitto 1.mo Trimestre 2021"; F = "z0bug.fy_2021" },
    @{ Row = 9;  A = "z0bug.move_8";  C = "####-04-05"; D = "Affitto 2.do Trimestre 2021"; F = "z0bug.fy_2021" },
    @{ Row = 10; A = "z0bug.move_9";  C = "####-07-05"; D = "Affitto 3.zo Trimestre 2021"; F = "z0bug.fy_2021" },
    @{ Row = 11; A = "z0bug.move_10"; C = "####-10-05"; D = "Affitto 4.to Trimestre 2021"; F = "z0bug.fy_2021" }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column C carries the "date" text and is formatted as plain text
    # (number format "@"), same as the existing rows above.
    $ws.Range("C$row").NumberFormat = "@"

    $ws.Range("A$row").Value = $r.A
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = "z0bug.misc"
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = "entry"
}
